$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "asd"
$ws.Range("A3").Value = "asd"
$ws.Range("A4").Value = "asd"
$ws.Range("A5").Value = "asd"
$ws.Range("A6").Value = "asd"
$ws.Range("A7").Value = "ads"

$ws.Range("A7").Select()
